# Applies the 2022-05-16 Fonds de solidarite data update:
# updates nombre_aides (C), nombre_entreprises (D), montant_total (E) for
# the rows impacted by the refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 88;  C = 71266;    E = 110295157 },
    @{ Row = 91;  C = 18857;    E = 75163013 },
    @{ Row = 97;  C = 22;       E = 354577 },
    @{ Row = 98;  C = 6295;     E = 19452171 },
    @{ Row = 100; C = 9339;     E = 23760393 },
    @{ Row = 121; C = 1306135;  E = 2274587611 },
    @{ Row = 122; C = 366;      E = 1162215 },
    @{ Row = 129; C = 633363;   E = 3427264164 },
    @{ Row = 132; C = 585656;   E = 3462321824 },
    @{ Row = 136; C = 26677;    E = 143561800 },
    @{ Row = 139; C = 76641;    E = 114132538 },
    @{ Row = 144; C = 25074;    E = 92430892 },
    @{ Row = 146; C = 7440;     E = 37684204 },
    @{ Row = 151; C = 39928;    E = 60363412 },
    @{ Row = 154; C = 18453;    E = 72760118 },
    @{ Row = 156; C = 12398;    E = 40047588 },
    @{ Row = 171; C = 95821;    E = 490672924 },
    @{ Row = 189; C = 100465;   E = 556064990 },
    @{ Row = 224; C = 39675;    E = 260923034 },
    @{ Row = 229; C = 612545;   D = 121244; E = 1040753491 },
    @{ Row = 240; C = 205901;   E = 1068739496 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = $u.C
    if ($u.ContainsKey('D')) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
